$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: BASIC APP CHANGES (sheet7) -- finish the row-4 test case and append
# the new rows 5-13 documenting the completed testing pass.
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("BASIC APP CHANGES")

# Row 4: swap the zip file used and record the actual test outcome columns.
$ws7.Cells.Item(4, 4).Value = "BASIC-APP2.ZIP"
$ws7.Cells.Item(4, 5).Value = "Input start and end date as 6April 2015, select filter as CCA"
$ws7.Cells.Item(4, 6).Value = "75% choir, 25% soccer"
$ws7.Cells.Item(4, 7).Value = "75% choir, 25% soccer"
$ws7.Cells.Item(4, 8).Value = "Pass"
$ws7.Cells.Item(4, 2).WrapText = $true
$ws7.Cells.Item(4, 3).WrapText = $true
$ws7.Cells.Item(4, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(4, 5).WrapText = $true
$ws7.Cells.Item(4, 6).WrapText = $true
$ws7.Cells.Item(4, 7).WrapText = $true
$ws7.Cells.Item(4, 8).WrapText = $true
$ws7.Rows.Item(4).RowHeight = 116

# Row 5
$ws7.Cells.Item(5, 1).Value = 2
$ws7.Cells.Item(5, 2).Value = "Change to demographics.csv & Breakdown by usage time category and demographics Report"
$ws7.Cells.Item(5, 3).Value = "User are able to sort by Gender accurately"
$ws7.Cells.Item(5, 4).Value = "BASIC-APP.ZIP"
$ws7.Cells.Item(5, 5).Value = "Input start and end date as 6April 2015, select filter as gender"
$ws7.Cells.Item(5, 6).Value = "50% Female, 50% Male, both mild"
$ws7.Cells.Item(5, 7).Value = "50% Female, 50% Male, both mild"
$ws7.Cells.Item(5, 8).Value = "Pass"
$ws7.Cells.Item(5, 2).WrapText = $true
$ws7.Cells.Item(5, 3).WrapText = $true
$ws7.Cells.Item(5, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(5, 5).WrapText = $true
$ws7.Cells.Item(5, 6).WrapText = $true
$ws7.Cells.Item(5, 7).WrapText = $true
$ws7.Cells.Item(5, 8).WrapText = $true
$ws7.Rows.Item(5).RowHeight = 116

# Row 6
$ws7.Cells.Item(6, 1).Value = 3
$ws7.Cells.Item(6, 2).Value = "Change to demographics.csv & Breakdown by usage time category and demographics Report"
$ws7.Cells.Item(6, 3).Value = "User are able to sort by School accurately"
$ws7.Cells.Item(6, 4).Value = "BASIC-APP.ZIP"
$ws7.Cells.Item(6, 5).Value = "Input start and end date as 6April 2015, select filter as School"
$ws7.Cells.Item(6, 6).Value = "1 mild user from econs, 1 from accountancy and 2 from sosci"
$ws7.Cells.Item(6, 7).Value = "1 mild user from econs, 1 from accountancy and 2 from sosci"
$ws7.Cells.Item(6, 8).Value = "Pass"
$ws7.Cells.Item(6, 2).WrapText = $true
$ws7.Cells.Item(6, 3).WrapText = $true
$ws7.Cells.Item(6, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(6, 5).WrapText = $true
$ws7.Cells.Item(6, 6).WrapText = $true
$ws7.Cells.Item(6, 7).WrapText = $true
$ws7.Cells.Item(6, 8).WrapText = $true
$ws7.Rows.Item(6).RowHeight = 130.5

# Row 7
$ws7.Cells.Item(7, 1).Value = 4
$ws7.Cells.Item(7, 2).Value = "Change to demographics.csv & Breakdown by usage time category and demographics Report"
$ws7.Cells.Item(7, 3).Value = "User are able to sort by Year accurately"
$ws7.Cells.Item(7, 4).Value = "BASIC-APP.ZIP"
$ws7.Cells.Item(7, 5).Value = "Input start and end date as 6April 2015, select filter as Year"
$ws7.Cells.Item(7, 6).Value = "1 from 2015, 2 from 2014, 1 from 2012"
$ws7.Cells.Item(7, 7).Value = "1 from 2015, 2 from 2014, 1 from 2012"
$ws7.Cells.Item(7, 8).Value = "Pass"
$ws7.Cells.Item(7, 2).WrapText = $true
$ws7.Cells.Item(7, 3).WrapText = $true
$ws7.Cells.Item(7, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(7, 5).WrapText = $true
$ws7.Cells.Item(7, 6).WrapText = $true
$ws7.Cells.Item(7, 7).WrapText = $true
$ws7.Cells.Item(7, 8).WrapText = $true
$ws7.Rows.Item(7).RowHeight = 116

# Row 8
$ws7.Cells.Item(8, 1).Value = 5
$ws7.Cells.Item(8, 2).Value = "Change to demographics.csv & Breakdown by usage time category and demographics Report"
$ws7.Cells.Item(8, 3).Value = "Test the filters, if I can randomly select a filter if I only have 1 demographics to sort out in mind"
$ws7.Cells.Item(8, 4).Value = "UI"
$ws7.Cells.Item(8, 5).Value = "Try out all filter"
$ws7.Cells.Item(8, 6).Value = "All filters can work"
$ws7.Cells.Item(8, 7).Value = "All filters can work"
$ws7.Cells.Item(8, 8).Value = "Pass"
$ws7.Cells.Item(8, 2).WrapText = $true
$ws7.Cells.Item(8, 3).WrapText = $true
$ws7.Cells.Item(8, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(8, 5).WrapText = $true
$ws7.Cells.Item(8, 6).WrapText = $true
$ws7.Cells.Item(8, 7).WrapText = $true
$ws7.Cells.Item(8, 8).WrapText = $true
$ws7.Rows.Item(8).RowHeight = 43.5

# Row 9
$ws7.Cells.Item(9, 1).Value = 6
$ws7.Cells.Item(9, 2).Value = "Change to demographics.csv & Breakdown by usage time category and demographics Report"
$ws7.Cells.Item(9, 3).Value = "Test if the calculation of Intense, Moderate and Mild work"
$ws7.Cells.Item(9, 4).Value = "BASIC-APP3.ZIP"
$ws7.Cells.Item(9, 5).Value = "Input start and end date as 6April 2015, select filter as gender"
$ws7.Cells.Item(9, 6).Value = "50% Female both mild, 50% Male with one normal and one mild "
$ws7.Cells.Item(9, 7).Value = "50% Female both mild, 50% Male with one normal and one mild "
$ws7.Cells.Item(9, 8).Value = "Pass"
$ws7.Cells.Item(9, 2).WrapText = $true
$ws7.Cells.Item(9, 3).WrapText = $true
$ws7.Cells.Item(9, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(9, 5).WrapText = $true
$ws7.Cells.Item(9, 6).WrapText = $true
$ws7.Cells.Item(9, 7).WrapText = $true
$ws7.Cells.Item(9, 8).WrapText = $true
$ws7.Rows.Item(9).RowHeight = 145

# Row 10
$ws7.Cells.Item(10, 1).Value = 7
$ws7.Cells.Item(10, 2).Value = "Change to demographics.csv & Breakdown by usage time category and demographics Report"
$ws7.Cells.Item(10, 3).Value = "Test if the average function work"
$ws7.Cells.Item(10, 4).Value = "BASIC-APP2.ZIP"
$ws7.Cells.Item(10, 5).Value = "Retry row 2-4, but this time date is 06/04/2015 to 07/04/2015"
$ws7.Cells.Item(10, 6).Value = "Result of all fields should be similar to row 2-4"
$ws7.Cells.Item(10, 7).Value = "Result of all fields same as what row 2-4 returns"
$ws7.Cells.Item(10, 8).Value = "Pass"
$ws7.Cells.Item(10, 2).WrapText = $true
$ws7.Cells.Item(10, 3).WrapText = $true
$ws7.Cells.Item(10, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(10, 5).WrapText = $true
$ws7.Cells.Item(10, 6).WrapText = $true
$ws7.Cells.Item(10, 7).WrapText = $true
$ws7.Cells.Item(10, 8).WrapText = $true
$ws7.Rows.Item(10).RowHeight = 130.5

# Row 11
$ws7.Cells.Item(11, 1).Value = 8
$ws7.Cells.Item(11, 2).Value = "Change to demographics.csv & Breakdown by usage time category and demographics Report"
$ws7.Cells.Item(11, 3).Value = "Check of valid start date"
$ws7.Cells.Item(11, 4).Value = "sets start date as 4-12-0001"
$ws7.Cells.Item(11, 5).Value = "UI"
$ws7.Cells.Item(11, 6).Value = """invalid startdate"""
$ws7.Cells.Item(11, 7).Value = """value must be 01/01/1970 and later"""
$ws7.Cells.Item(11, 8).Value = "Pass"
$ws7.Cells.Item(11, 2).WrapText = $true
$ws7.Cells.Item(11, 3).WrapText = $true
$ws7.Cells.Item(11, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(11, 4).HorizontalAlignment = -4108
$ws7.Cells.Item(11, 5).HorizontalAlignment = -4108
$ws7.Cells.Item(11, 6).HorizontalAlignment = -4108
$ws7.Cells.Item(11, 7).WrapText = $true
$ws7.Cells.Item(11, 8).WrapText = $true
$ws7.Rows.Item(11).RowHeight = 43.5

# Row 12
$ws7.Cells.Item(12, 1).Value = 9
$ws7.Cells.Item(12, 2).Value = "Change to demographics.csv & Breakdown by usage time category and demographics Report"
$ws7.Cells.Item(12, 3).Value = "Check for valid end date"
$ws7.Cells.Item(12, 4).Value = "sets end date as 3-12-0001"
$ws7.Cells.Item(12, 5).Value = "UI"
$ws7.Cells.Item(12, 6).Value = """invalid enddate"""
$ws7.Cells.Item(12, 7).Value = """value must be 01/01/1970 and later"""
$ws7.Cells.Item(12, 8).Value = "Pass"
$ws7.Cells.Item(12, 2).WrapText = $true
$ws7.Cells.Item(12, 3).WrapText = $true
$ws7.Cells.Item(12, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(12, 4).HorizontalAlignment = -4108
$ws7.Cells.Item(12, 5).HorizontalAlignment = -4108
$ws7.Cells.Item(12, 6).HorizontalAlignment = -4108
$ws7.Cells.Item(12, 7).WrapText = $true
$ws7.Cells.Item(12, 8).WrapText = $true
$ws7.Rows.Item(12).RowHeight = 43.5

# Row 13
$ws7.Cells.Item(13, 1).Value = 10
$ws7.Cells.Item(13, 2).Value = "Change to demographics.csv & Breakdown by usage time category and demographics Report"
$ws7.Cells.Item(13, 3).Value = "Check for valid start date before end date"
$ws7.Cells.Item(13, 4).Value = "start date : 07/04/2015" + [char]10 + "end date: 06/04/2015"
$ws7.Cells.Item(13, 5).Value = "UI"
$ws7.Cells.Item(13, 6).Value = "end date must be after start date"
$ws7.Cells.Item(13, 7).Value = "end date must be after start date"
$ws7.Cells.Item(13, 8).Value = "Pass"
$ws7.Cells.Item(13, 2).WrapText = $true
$ws7.Cells.Item(13, 3).WrapText = $true
$ws7.Cells.Item(13, 3).HorizontalAlignment = -4108
$ws7.Cells.Item(13, 4).WrapText = $true
$ws7.Cells.Item(13, 4).HorizontalAlignment = -4108
$ws7.Cells.Item(13, 5).HorizontalAlignment = -4108
$ws7.Cells.Item(13, 6).WrapText = $true
$ws7.Cells.Item(13, 6).HorizontalAlignment = -4108
$ws7.Cells.Item(13, 7).WrapText = $true
$ws7.Cells.Item(13, 7).HorizontalAlignment = -4108
$ws7.Cells.Item(13, 8).WrapText = $true
$ws7.Rows.Item(13).RowHeight = 43.5

# Column widths for the newly-used columns E:G
$ws7.Columns.Item(5).ColumnWidth = 7.5
$ws7.Columns.Item(6).ColumnWidth = 15.2
$ws7.Columns.Item(7).ColumnWidth = 15.2

# ---------------------------------------------------------------------------
# Sheet: SOCIAL ACTIVENESS (sheet3) -- move the selection onto the results
# table that was being reviewed.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("SOCIAL ACTIVENESS")
$ws3.Activate()
$ws3.Range("C4:F6").Select()

# ---------------------------------------------------------------------------
# Finally land on BASIC APP CHANGES (this also clears tabSelected from
# BOOTSTRAP CHANGES and makes this the workbook's active tab).
# ---------------------------------------------------------------------------
$ws7.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws7.Range("A14").Select()
